$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.214.82"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "1.661.62"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.61"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5233"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2668"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06316"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.05"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07727"
$ws.Range("D12").Value = "1.656.04"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.427"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "1.889.14"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5471"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").Value = "0.0₅8180"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.86"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "26.249.17"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.663"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.96"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.15"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.077"
$ws.Range("E23").Value = "  -4.20%  "
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.02"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1237"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.14"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05997"
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.659"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.311"
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9799"
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.785"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.413"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5870"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01590"
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.946"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8613"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "1.033.12"
$ws.Range("E43").Value = "  -3.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.58"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("D45").Value = "1.803.73"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.078"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.466"
$ws.Range("E51").Value = "  -0.32%  "
